# Apply BOM updates: added pull-up resistor to NRST line, changed I2C and PWM functions mapping
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capacitor rows (0.1uF / 1uF groups) ---
# Row 3 becomes the 0.1uF (100nF) group, row 4 becomes the 1uF group
$ws.Range("A3").Value = "C2, C5, C6, C7, C8, C9, C12"
$ws.Range("B3").Value = "0.1uF ± 10%, 25V, X7R, SMD 0603"
$ws.Range("C3").Value = "100nF"
$ws.Range("D3").Value = 7

$ws.Range("A4").Value = "C3, C4, C11"
$ws.Range("B4").Value = "1uF ± 10%, 16V, X7R, SMD 0603"
$ws.Range("C4").Value = "1uF"
$ws.Range("D4").Value = 3

# --- Resistor renumbering: new pull-up resistor R6 added to NRST line ---
$ws.Range("A15").Value = "R5, R6, R10, R11"
$ws.Range("D15").Value = 4

$ws.Range("A16").Value = "R7, R8"

$ws.Range("A17").Value = "R9"

# --- Connector function mapping changes (I2C / PWM) ---
$ws.Range("B20").Value = "PLS-5"
$ws.Range("C20").Value = "PROG"

$ws.Range("B21").Value = "PLD-4"
$ws.Range("C21").Value = "BOOT"

$ws.Range("C26").Value = "AXIS-A4"
$ws.Range("C27").Value = "SPI/A5"
$ws.Range("C28").Value = "SPI/A6"
$ws.Range("C29").Value = "SHIFT REG/A7"

# Row heights follow the wrapped-text content (row3 now has the long 0.1uF
# description that wraps to two lines, row4 now has the shorter 1uF
# description that fits a single line)
$ws.Rows(4).AutoFit()
$ws.Rows(3).RowHeight = 30
